$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 416.66666
$ws.Range("I6").Value = 400
$ws.Range("K6").Value = 1200
$ws.Range("M6").Value = -1088
$ws.Range("H28").Value = 776.8570999999999
$ws.Range("I28").Value = 767.6
$ws.Range("K28").Value = 767.6
$ws.Range("M28").Value = -282.6
$ws.Range("H86").Value = 3664.0908
$ws.Range("I86").Value = 2252.125
$ws.Range("J86").Value = 4993
$ws.Range("K86").Value = 2252.125
$ws.Range("L86").Value = 4993
$ws.Range("M86").Value = -1129.125
$ws.Range("N86").Value = -7239
$ws.Range("H89").Value = 3664.0908
$ws.Range("I89").Value = 2252.125
$ws.Range("J89").Value = 4993
$ws.Range("K89").Value = 11260.625
$ws.Range("L89").Value = 24965
$ws.Range("M89").Value = -5644.625
$ws.Range("N89").Value = -36197
$ws.Range("H116").Value = 120476.89
$ws.Range("I116").Value = 178215.33
$ws.Range("K116").Value = 178215.33
$ws.Range("M116").Value = -174773.33

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3344.38
$ws.Range("I32").Value = 2907.9463
$ws.Range("J32").Value = 9142.714
$ws.Range("K32").Value = 2907.9463
$ws.Range("L32").Value = 9142.714
$ws.Range("M32").Value = -2620.9463
$ws.Range("N32").Value = -9716.714
$ws.Range("H45").Value = 1221.2916
$ws.Range("I45").Value = 1142.4615
$ws.Range("J45").Value = 1314.4546
$ws.Range("K45").Value = 1142.4615
$ws.Range("L45").Value = 1314.4546
$ws.Range("M45").Value = -765.4614999999999
$ws.Range("N45").Value = -2068.4546

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1864.1111
$ws.Range("I107").Value = 1996.0667
$ws.Range("K107").Value = 1996.0667
$ws.Range("M107").Value = -76.06670000000008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2345.889
$ws.Range("I58").Value = 1848
$ws.Range("J58").Value = 2627.3044
$ws.Range("K58").Value = 1848
$ws.Range("L58").Value = 2627.3044
$ws.Range("M58").Value = -1645
$ws.Range("N58").Value = -3033.3044
$ws.Range("H122").Value = 2501.7144
$ws.Range("I122").Value = 2585.3333
$ws.Range("K122").Value = 7755.999899999999
$ws.Range("M122").Value = -5305.999899999999
$ws.Range("H132").Value = 1717.1805
$ws.Range("I132").Value = 1000.4808
$ws.Range("J132").Value = 3580.6
$ws.Range("K132").Value = 3001.4424
$ws.Range("L132").Value = 10741.8
$ws.Range("M132").Value = -471.4423999999999
$ws.Range("N132").Value = -15801.8
$ws.Range("H134").Value = 1358.3969
$ws.Range("I134").Value = 1018.6889
$ws.Range("J134").Value = 2207.6667
$ws.Range("K134").Value = 3056.0667
$ws.Range("L134").Value = 6623.000100000001
$ws.Range("M134").Value = -521.0666999999999
$ws.Range("N134").Value = -11693.0001
$ws.Range("H136").Value = 2345.889
$ws.Range("I136").Value = 1848
$ws.Range("J136").Value = 2627.3044
$ws.Range("K136").Value = 5544
$ws.Range("L136").Value = 7881.9132
$ws.Range("M136").Value = -2994
$ws.Range("N136").Value = -12981.9132

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H12").Value = 74.458336
$ws.Range("J12").Value = 95.29412000000001
$ws.Range("L12").Value = 285.88236
$ws.Range("N12").Value = -631.8823600000001
$ws.Range("H16").Value = 614.2857
$ws.Range("I16").Value = 460
$ws.Range("J16").Value = 1000
$ws.Range("K16").Value = 1380
$ws.Range("L16").Value = 3000
$ws.Range("M16").Value = -1207
$ws.Range("N16").Value = -3346
$ws.Range("H38").Value = 59.42857
$ws.Range("I38").Value = 69.333336
$ws.Range("J38").Value = 41.6
$ws.Range("K38").Value = 208.000008
$ws.Range("L38").Value = 124.8
$ws.Range("M38").Value = 138.999992
$ws.Range("N38").Value = -818.8
$ws.Range("H96").Value = 3633.3333
$ws.Range("J96").Value = 3633.3333
$ws.Range("L96").Value = 10899.9999
$ws.Range("N96").Value = -15017.9999
$ws.Range("H131").Value = 1429.1177
$ws.Range("I131").Value = 599.3158
$ws.Range("J131").Value = 1921.8125
$ws.Range("K131").Value = 1797.9474
$ws.Range("L131").Value = 5765.4375
$ws.Range("M131").Value = 3242.0526
$ws.Range("N131").Value = -15845.4375

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H82").Value = 0
$ws.Range("J82").Value = 0
$ws.Range("L82").Value = 0
$ws.Range("H85").Value = 0
$ws.Range("J85").Value = 0
$ws.Range("L85").Value = 0
$ws.Range("H102").Value = 2037.3103
$ws.Range("I102").Value = 2024.3572
$ws.Range("K102").Value = 2024.3572
$ws.Range("M102").Value = -402.3571999999999
$ws.Range("H132").Value = 2916.4807
$ws.Range("I132").Value = 2765.725
$ws.Range("J132").Value = 3419
$ws.Range("K132").Value = 8297.174999999999
$ws.Range("L132").Value = 10257
$ws.Range("M132").Value = -5767.174999999999
$ws.Range("N132").Value = -15317
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4140.6
$ws.Range("I16").Value = 3740.8333
$ws.Range("J16").Value = 5739.6665
$ws.Range("K16").Value = 3740.8333
$ws.Range("L16").Value = 5739.6665
$ws.Range("M16").Value = -3570.8333
$ws.Range("N16").Value = -6079.6665
$ws.Range("H32").Value = 4766.6665
$ws.Range("I32").Value = 2900
$ws.Range("K32").Value = 2900
$ws.Range("M32").Value = -2583
$ws.Range("H46").Value = 770
$ws.Range("I46").Value = 695.55554
$ws.Range("J46").Value = 993.3333
$ws.Range("K46").Value = 695.55554
$ws.Range("L46").Value = 993.3333
$ws.Range("M46").Value = -507.55554
$ws.Range("N46").Value = -1369.3333
$ws.Range("H74").Value = 28649.5
$ws.Range("J74").Value = 28649.5
$ws.Range("L74").Value = 28649.5
$ws.Range("N74").Value = -30645.5
$ws.Range("H77").Value = 28649.5
$ws.Range("J77").Value = 28649.5
$ws.Range("L77").Value = 85948.5
$ws.Range("N77").Value = -95932.5
$ws.Range("H122").Value = 3393.1904
$ws.Range("I122").Value = 3310.4666
$ws.Range("J122").Value = 3600
$ws.Range("K122").Value = 9931.399800000001
$ws.Range("L122").Value = 10800
$ws.Range("M122").Value = -7481.399800000001
$ws.Range("N122").Value = -15700
$ws.Range("H132").Value = 11622.591
$ws.Range("I132").Value = 3722.111
$ws.Range("J132").Value = 17092.154
$ws.Range("K132").Value = 11166.333
$ws.Range("L132").Value = 51276.462
$ws.Range("M132").Value = -8636.332999999999
$ws.Range("N132").Value = -56336.462

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 736.53845
$ws.Range("I107").Value = 898.4
$ws.Range("J107").Value = 635.375
$ws.Range("K107").Value = 2695.2
$ws.Range("L107").Value = 1906.125
$ws.Range("M107").Value = -775.1999999999998
$ws.Range("N107").Value = -5746.125
